$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$c1 = $ws2.Range("C1")
$c1.Style = "Normal"
$c1.BorderAround(1)
$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142

$d1 = $ws2.Range("D1")
$d1.Style = "Normal"
$d1.BorderAround(1)
$d1.Borders.Item(7).LineStyle = -4142

$f1 = $ws2.Range("F1")
$f1.Style = "Normal"
$f1.Borders.Item(8).LineStyle = 1
$f1.Borders.Item(9).LineStyle = 1

$g1 = $ws2.Range("G1")
$g1.Style = "Normal"
$g1.Borders.Item(8).LineStyle = 1
$g1.Borders.Item(9).LineStyle = 1
$g1.Borders.Item(10).LineStyle = 1

$ws1 = $wb.Worksheets.Item(1)

$a = $ws1.Range("C1")
$a.Style = "Normal"
$a.Borders.Item(8).LineStyle = 1
$a.Borders.Item(9).LineStyle = 1

$b = $ws1.Range("D1")
$b.Style = "Normal"
$b.Borders.Item(8).LineStyle = 1
$b.Borders.Item(9).LineStyle = 1
$b.Borders.Item(10).LineStyle = 1
